# Estudio previo contrato prestacion de servicios - targeted edit
#
# The only substantive textual change introduced by the tracked revision is
# the addition of the phrase " mismo valor del inicio" right after the
# "(______ ($____)auto" placeholder that precedes "incluido IVA, impuestos o
# descuentos." (the remainder of the underlying diff is purely cosmetic
# <w:proofErr/> bookkeeping and <w:lastRenderedPageBreak/> repagination noise
# that Word's proofing engine / layout engine regenerates on its own and
# carries no visible-text or object-model-observable effect).
#
# We locate the unique run of text "suma de: ______ ($____)auto" with
# Find.Execute (no replacement), collapse the resulting range to its end
# point (i.e. right after "auto"), and insert the new text there with
# InsertAfter so the existing runs/formatting around it are left completely
# untouched and only the new text gets spliced in between "auto" and the
# following " incluido IVA, impuestos o descuentos." run.

$d = $word.ActiveDocument

$search = $d.Content
$found = $search.Find.Execute(
    "suma de: ______ (`$____)auto",  # FindText
    $true,                           # MatchCase
    $false,                          # MatchWholeWord
    $false,                          # MatchWildcards
    $false,                          # MatchSoundsLike
    $false,                          # MatchAllWordForms
    $true,                           # Forward
    1,                               # Wrap (wdFindContinue)
    $false,                          # Format
    "",                              # ReplaceWith
    0                                # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not locate the target placeholder text to edit."
}

$insertionPoint = $d.Range($search.End, $search.End)
$insertionPoint.InsertAfter(" mismo valor del inicio")
